# Auto-generated edit script applying scheduled market-data refresh values
# to the Seraph_Profits workbook, per the provided OOXML diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 162.33333
$ws.Range("I33").Value = 152.5
$ws.Range("K33").Value = 152.5
$ws.Range("M33").Value = 76.5

# Row 62
$ws.Range("H62").Value = 7272
$ws.Range("J62").Value = 9049.875
$ws.Range("L62").Value = 9049.875
$ws.Range("N62").Value = -10297.875

# Row 65
$ws.Range("H65").Value = 7272
$ws.Range("J65").Value = 9049.875
$ws.Range("L65").Value = 45249.375
$ws.Range("N65").Value = -51489.375

# Row 69
$ws.Range("H69").Value = 12000
$ws.Range("J69").Value = 12000
$ws.Range("L69").Value = 36000
$ws.Range("N69").Value = -37748

# Row 72
$ws.Range("H72").Value = 12000
$ws.Range("J72").Value = 12000
$ws.Range("L72").Value = 108000
$ws.Range("N72").Value = -116736

# Row 76
$ws.Range("H76").Value = 7287.5
$ws.Range("I76").Value = 4400
$ws.Range("K76").Value = 4400
$ws.Range("M76").Value = -4085

# Row 79
$ws.Range("H79").Value = 7287.5
$ws.Range("I79").Value = 4400
$ws.Range("K79").Value = 4400
$ws.Range("M79").Value = -3308

# Row 99
$ws.Range("H99").Value = 289.14285

# Row 118
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()

# Row 138
$ws.Range("H138").Value = 3265.5264
$ws.Range("I138").Value = 1653.6364
$ws.Range("J138").Value = 3538.3076
$ws.Range("K138").Value = 4960.9092
$ws.Range("L138").Value = 10614.9228
$ws.Range("M138").Value = 179.0907999999999
$ws.Range("N138").Value = -20894.9228


$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 4472.8
$ws.Range("I61").Value = 4472.8
$ws.Range("K61").Value = 4472.8
$ws.Range("M61").Value = -4260.8

# Row 63
$ws.Range("H63").Value = 5756.5713
$ws.Range("I63").Value = 3216.5
$ws.Range("K63").Value = 3216.5
$ws.Range("M63").Value = -2530.5

# Row 66
$ws.Range("H66").Value = 5756.5713
$ws.Range("I66").Value = 3216.5
$ws.Range("K66").Value = 16082.5
$ws.Range("M66").Value = -12650.5

# Row 95
$ws.Range("H95").Value = 50367.8
$ws.Range("J95").Value = 50367.8
$ws.Range("L95").Value = 50367.8
$ws.Range("N95").Value = -55859.8

# Row 122
$ws.Range("H122").Value = 2352.5417
$ws.Range("I122").Value = 1527.1177
$ws.Range("J122").Value = 4357.143
$ws.Range("K122").Value = 4581.3531
$ws.Range("L122").Value = 13071.429
$ws.Range("M122").Value = -2131.3531
$ws.Range("N122").Value = -17971.429

# Row 132
$ws.Range("H132").Value = 2988.2144
$ws.Range("I132").Value = 2530.5454
$ws.Range("J132").Value = 4666.3335
$ws.Range("K132").Value = 7591.6362
$ws.Range("L132").Value = 13999.0005
$ws.Range("M132").Value = -5061.6362
$ws.Range("N132").Value = -19059.0005

# Row 136
$ws.Range("H136").Value = 4472.8
$ws.Range("I136").Value = 4472.8
$ws.Range("K136").Value = 13418.4
$ws.Range("M136").Value = -10868.4


$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 938.9231
$ws.Range("I94").Value = 991.3913
$ws.Range("K94").Value = 991.3913
$ws.Range("M94").Value = -540.3913

# Row 99
$ws.Range("H99").Value = 3980
$ws.Range("I99").Value = 3976.8333
$ws.Range("K99").Value = 3976.8333
$ws.Range("M99").Value = -2478.8333


$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4338.2666
$ws.Range("I31").Value = 2717
$ws.Range("K31").Value = 2717
$ws.Range("M31").Value = -2422

# Row 34
$ws.Range("H34").Value = 4338.2666
$ws.Range("I34").Value = 2717
$ws.Range("K34").Value = 2717
$ws.Range("M34").Value = -2515

# Row 76
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

# Row 79
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

# Row 86
$ws.Range("H86").Value = 8519.4
$ws.Range("I86").Value = 7232.3335
$ws.Range("K86").Value = 7232.3335
$ws.Range("M86").Value = -6109.3335

# Row 89
$ws.Range("H89").Value = 8519.4
$ws.Range("I89").Value = 7232.3335
$ws.Range("K89").Value = 36161.6675
$ws.Range("M89").Value = -30545.6675


$ws = $wb.Worksheets.Item("CUL")
# Row 24
$ws.Range("H24").Value = 392.2857
$ws.Range("J24").Value = 416
$ws.Range("L24").Value = 1248
$ws.Range("N24").Value = -1708

# Row 101
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

# Row 117
$ws.Range("H117").Value = 1272.4546
$ws.Range("I117").Value = 251.875
$ws.Range("K117").Value = 755.625
$ws.Range("M117").Value = 2686.375

# Row 131
$ws.Range("H131").Value = 1377.3334
$ws.Range("I131").Value = 1029.7142
$ws.Range("J131").Value = 1499
$ws.Range("K131").Value = 3089.1426
$ws.Range("L131").Value = 4497
$ws.Range("M131").Value = 1950.8574
$ws.Range("N131").Value = -14577


$ws = $wb.Worksheets.Item("GSM")
# Row 20
$ws.Range("H20").Value = 45872.668
$ws.Range("J20").Value = 59644.547
$ws.Range("L20").Value = 59644.547
$ws.Range("N20").Value = -60134.547


$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 359.16666
$ws.Range("I55").Value = 380.14285
$ws.Range("J55").Value = 285.75
$ws.Range("K55").Value = 380.14285
$ws.Range("L55").Value = 285.75
$ws.Range("M55").Value = -207.14285
$ws.Range("N55").Value = -631.75

# Row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

# Row 136
$ws.Range("H136").Value = 4000
$ws.Range("I136").Value = 4000
$ws.Range("K136").Value = 12000
$ws.Range("M136").Value = -9450


$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Range("H45").Value = 19286.1
$ws.Range("J45").Value = 17745.25
$ws.Range("L45").Value = 17745.25
$ws.Range("N45").Value = -18727.25

# Row 81
$ws.Range("H81").Value = 2508.6667
$ws.Range("I81").Value = 1972.9231
$ws.Range("K81").Value = 3945.8462
$ws.Range("M81").Value = -2884.8462

# Row 84
$ws.Range("H84").Value = 2508.6667
$ws.Range("I84").Value = 1972.9231
$ws.Range("K84").Value = 19729.231
$ws.Range("M84").Value = -14425.231

# Row 122
$ws.Range("H122").Value = 2998.4
$ws.Range("I122").Value = 3548.0833
$ws.Range("K122").Value = 10644.2499
$ws.Range("M122").Value = -8194.249899999999

# Row 132
$ws.Range("H132").Value = 1817.1666
$ws.Range("I132").Value = 1634.3334
$ws.Range("K132").Value = 4903.0002
$ws.Range("M132").Value = -2373.0002

# Row 136
$ws.Range("H136").Value = 8193
$ws.Range("I136").Value = 8990.333000000001
$ws.Range("J136").Value = 6997
$ws.Range("K136").Value = 26970.999
$ws.Range("L136").Value = 20991
$ws.Range("M136").Value = -24420.999
$ws.Range("N136").Value = -26091

